$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1889
$ws1.Range("F4").Value = 884
$ws1.Range("F5").Value = 772
$ws1.Range("F6").Value = 13343
$ws1.Range("F7").Value = 13214
$ws1.Range("F8").Value = 1020
$ws1.Range("F9").Value = 775
$ws1.Range("F10").Value = 23
$ws1.Range("F11").Value = 563
$ws1.Range("F13").Value = 681
$ws1.Range("F20").Value = 263
$ws1.Range("F22").Value = 424
$ws1.Range("F24").Value = 20

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 95
$ws2.Range("F3").Value = 35
$ws2.Range("F7").Value = 126

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 47

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1889
$ws4.Range("F5").Value = 884
$ws4.Range("F6").Value = 95
$ws4.Range("F7").Value = 772
$ws4.Range("F8").Value = 13343
$ws4.Range("F9").Value = 13214
$ws4.Range("F10").Value = 1020
$ws4.Range("F11").Value = 775
$ws4.Range("F12").Value = 23
$ws4.Range("F13").Value = 563
$ws4.Range("F15").Value = 681
$ws4.Range("F16").Value = 35
$ws4.Range("F25").Value = 47
$ws4.Range("F27").Value = 263
$ws4.Range("F29").Value = 424
$ws4.Range("F31").Value = 126
$ws4.Range("F33").Value = 20

$wb.Save()
